# Setting up inital code
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Palendrome")

# Replace the placeholder "<Your Name>" (cell C6) with the author's name.
# Cells C11, C12, C15, C16 hold formulas "=$C$6" so they pick this up
# automatically on recalculation.
$ws.Range("C6").Value = "Isaiah"

# Update the visible scroll position / selection on the sheet
# (mirrors the author scrolling down and selecting G15 before saving).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G15").Select()
